$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3550, 4360, 4729, 4729, 4807, 4812, 4890, 4890, 4890, 4890, 4890, 4890, 5053, 5053)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
